# Automatic update of files.
# Update the "Förändrad" (Changed) column (C) date from 2023-09-05 (45174)
# to 2023-09-06 (45175) for all data rows (rows 2 through 89).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 89
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
